$wb = $excel.ActiveWorkbook

# "Table-3.1" sheet holds the answers for the top-3 English speaking
# countries question. Fill in the actual solution values.
$ws = $wb.Worksheets.Item("Table-3.1")

$ws.Cells.Item(5, 3).Value = "United States (USA)"
$ws.Cells.Item(6, 3).Value = "United Kingdom (GBR)"
$ws.Cells.Item(7, 3).Value = "India (IND)"
